$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update item name labels in column A (rows 2-13) ---
# "full box of vials" (NO) was dropped from the inventory and replaced by
# "packaging of glass vials with penicillin G" (NO); "production of gloves" (MY) was
# renamed to "production of a pair of gloves" (MY). The full set of 12 item labels is
# written here in its new alphabetical order (A2:A13).
#
# Note: each label starts/ends with a literal single-quote character. Excel/COM treats a
# *leading* single quote in a .Value assignment as the 'force text' prefix marker and won't
# store it, so we double the leading quote (standard escape) to get a literal one back.
$ws.Range("A2").Value = "''market for sodium chlorate, powder' (kilogram, RER, None)"
$ws.Range("A3").Value = "''market for water, ultrapure' (kilogram, RER, None)"
$ws.Range("A4").Value = "''medical connector' (unit, GLO, None)"
$ws.Range("A5").Value = "''packaging of glass vials with penicillin G' (unit, NO, None)"
$ws.Range("A6").Value = "''packed box of penicillin' (unit, SE, None)"
$ws.Range("A7").Value = "''production of IV sets' (unit, RER, None)"
$ws.Range("A8").Value = "''production of a pair of gloves' (unit, MY, None)"
$ws.Range("A9").Value = "''production of alchohol wipes' (unit, DK, None)"
$ws.Range("A10").Value = "''stopcock' (unit, GLO, None)"
$ws.Range("A11").Value = "''treatment of hazardous waste, hazardous waste incineration' (kilogram, CH, None)"
$ws.Range("A12").Value = "''treatment of hazardous waste, hazardous waste incineration, with energy recovery' (kilogram, CH, None)"
$ws.Range("A13").Value = "''waste packaging paper, Recycled Content cut-off' (kilogram, GLO, None)"

# --- Update recomputed Monte-Carlo LCIA result values for rows 2-12 (row 13 is unchanged, all zeros) ---
# row 2
$ws.Range("B2").Value = [double]"0.008129370695716244"
$ws.Range("C2").Value = [double]"2.265880850336156"
$ws.Range("D2").Value = [double]"0.1796234922074676"
$ws.Range("E2").Value = [double]"0.2365080070114341"
$ws.Range("F2").Value = [double]"8.818108620020977"
$ws.Range("G2").Value = [double]"0.6072375609908904"
$ws.Range("H2").Value = [double]"0.001915808677453043"
$ws.Range("I2").Value = [double]"0.0002082153765173074"
$ws.Range("J2").Value = [double]"0.4312754708958648"
$ws.Range("K2").Value = [double]"4.10740454361543"
$ws.Range("L2").Value = [double]"1.190368312945629"
$ws.Range("M2").Value = [double]"0.07253609535142454"
$ws.Range("N2").Value = [double]"0.1359272064653131"
$ws.Range("O2").Value = [double]"1.063829392935665E-06"
$ws.Range("P2").Value = [double]"0.003317429232728206"
$ws.Range("Q2").Value = [double]"0.004422416867742256"
$ws.Range("R2").Value = [double]"0.004567556166330486"
$ws.Range("S2").Value = [double]"0.1200064134930158"
$ws.Range("T2").Value = [double]"1.083474859330304E-08"
$ws.Range("U2").Value = [double]"6.836226173468036E-06"
$ws.Range("V2").Value = [double]"0.1609293286047035"

# row 3
$ws.Range("B3").Value = [double]"8.49588496423555E-06"
$ws.Range("C3").Value = [double]"0.002608180301199517"
$ws.Range("D3").Value = [double]"0.0001124775847746949"
$ws.Range("E3").Value = [double]"0.0001498963660161549"
$ws.Range("F3").Value = [double]"0.005535752540869971"
$ws.Range("G3").Value = [double]"0.0006726734099012684"
$ws.Range("H3").Value = [double]"4.673999887934856E-05"
$ws.Range("I3").Value = [double]"6.403698784649995E-06"
$ws.Range("J3").Value = [double]"0.0005457582357520327"
$ws.Range("K3").Value = [double]"0.003018374494305449"
$ws.Range("L3").Value = [double]"0.001323482351792715"
$ws.Range("M3").Value = [double]"7.170881597748767E-05"
$ws.Range("N3").Value = [double]"6.976643490977473E-05"
$ws.Range("O3").Value = [double]"1.278067859528642E-09"
$ws.Range("P3").Value = [double]"3.521594376198354E-06"
$ws.Range("Q3").Value = [double]"4.721176863451457E-06"
$ws.Range("R3").Value = [double]"4.894953070982019E-06"
$ws.Range("S3").Value = [double]"0.0001057977629858891"
$ws.Range("T3").Value = [double]"4.184198443453184E-11"
$ws.Range("U3").Value = [double]"7.38379953345289E-09"
$ws.Range("V3").Value = [double]"0.0001580186552631967"

# row 4
$ws.Range("B4").Value = [double]"0.0005473478976743189"
$ws.Range("C4").Value = [double]"0.1114095287569723"
$ws.Range("D4").Value = [double]"0.02532106321273271"
$ws.Range("E4").Value = [double]"0.03746388221697419"
$ws.Range("F4").Value = [double]"0.9136830596574729"
$ws.Range("G4").Value = [double]"0.03223830631240234"
$ws.Range("H4").Value = [double]"6.307477853840808E-05"
$ws.Range("I4").Value = [double]"3.560176024460209E-06"
$ws.Range("J4").Value = [double]"0.01770091745324896"
$ws.Range("K4").Value = [double]"0.4456233253826027"
$ws.Range("L4").Value = [double]"0.004682702986805968"
$ws.Range("M4").Value = [double]"0.002402292253247127"
$ws.Range("N4").Value = [double]"0.003659590438655223"
$ws.Range("O4").Value = [double]"3.057861360296833E-08"
$ws.Range("P4").Value = [double]"0.0002361993781425839"
$ws.Range("Q4").Value = [double]"0.0003303129912550758"
$ws.Range("R4").Value = [double]"0.0003394339756427528"
$ws.Range("S4").Value = [double]"0.001089369841345707"
$ws.Range("T4").Value = [double]"5.673200565664793E-10"
$ws.Range("U4").Value = [double]"4.149488744705344E-07"
$ws.Range("V4").Value = [double]"0.009560071295497305"

# row 5
$ws.Range("B5").Value = [double]"0.003648298768081818"
$ws.Range("C5").Value = [double]"1.54909743418241"
$ws.Range("D5").Value = [double]"0.08621820806055434"
$ws.Range("E5").Value = [double]"0.1116120252805082"
$ws.Range("F5").Value = [double]"5.744517080799724"
$ws.Range("G5").Value = [double]"0.3785024802090458"
$ws.Range("H5").Value = [double]"0.0004711442001047956"
$ws.Range("I5").Value = [double]"0.0002836216033925794"
$ws.Range("J5").Value = [double]"0.063610407067083"
$ws.Range("K5").Value = [double]"1.321808912741056"
$ws.Range("L5").Value = [double]"0.2059085253813897"
$ws.Range("M5").Value = [double]"0.160747582042501"
$ws.Range("N5").Value = [double]"0.0284071520058528"
$ws.Range("O5").Value = [double]"1.045237715597826E-06"
$ws.Range("P5").Value = [double]"0.00143601970984707"
$ws.Range("Q5").Value = [double]"0.002600770093677633"
$ws.Range("R5").Value = [double]"0.00274311417134797"
$ws.Range("S5").Value = [double]"0.01952532260137308"
$ws.Range("T5").Value = [double]"7.344565491247429E-09"
$ws.Range("U5").Value = [double]"2.900468488545762E-06"
$ws.Range("V5").Value = [double]"0.1277338077318099"

# row 6
$ws.Range("B6").Value = [double]"0.002818221896174092"
$ws.Range("C6").Value = [double]"0.8251822069908367"
$ws.Range("D6").Value = [double]"0.04866138524179905"
$ws.Range("E6").Value = [double]"0.06077206530793378"
$ws.Range("F6").Value = [double]"3.888002494007186"
$ws.Range("G6").Value = [double]"0.2028633227751125"
$ws.Range("H6").Value = [double]"0.0003708441742537628"
$ws.Range("I6").Value = [double]"0.00021190091745003"
$ws.Range("J6").Value = [double]"0.04389101566798859"
$ws.Range("K6").Value = [double]"0.8864556607336793"
$ws.Range("L6").Value = [double]"0.1589854449635562"
$ws.Range("M6").Value = [double]"0.1147074106982291"
$ws.Range("N6").Value = [double]"0.02042903660015439"
$ws.Range("O6").Value = [double]"1.141032754490195E-06"
$ws.Range("P6").Value = [double]"0.001076236998287098"
$ws.Range("Q6").Value = [double]"0.001804046841105249"
$ws.Range("R6").Value = [double]"0.001881717175785557"
$ws.Range("S6").Value = [double]"0.02223195810650969"
$ws.Range("T6").Value = [double]"4.502260232904917E-09"
$ws.Range("U6").Value = [double]"1.84285611755583E-06"
$ws.Range("V6").Value = [double]"0.06425444695674611"

# row 7
$ws.Range("B7").Value = [double]"0.0005781492380916243"
$ws.Range("C7").Value = [double]"0.2026360647084039"
$ws.Range("D7").Value = [double]"0.006015986690459076"
$ws.Range("E7").Value = [double]"0.007915604025959396"
$ws.Range("F7").Value = [double]"0.5583825032869445"
$ws.Range("G7").Value = [double]"0.09515669650584613"
$ws.Range("H7").Value = [double]"5.07039174324563E-05"
$ws.Range("I7").Value = [double]"6.930443845205262E-06"
$ws.Range("J7").Value = [double]"0.009045009523788724"
$ws.Range("K7").Value = [double]"0.1331272485073503"
$ws.Range("L7").Value = [double]"0.01183923018480569"
$ws.Range("M7").Value = [double]"0.006223448989809889"
$ws.Range("N7").Value = [double]"0.00488558513888651"
$ws.Range("O7").Value = [double]"6.109224109761774E-08"
$ws.Range("P7").Value = [double]"0.0002646728661219914"
$ws.Range("Q7").Value = [double]"0.0004351791418714537"
$ws.Range("R7").Value = [double]"0.0004624809644617244"
$ws.Range("S7").Value = [double]"0.002271537471322852"
$ws.Range("T7").Value = [double]"8.503408545241858E-10"
$ws.Range("U7").Value = [double]"4.202994949806303E-07"
$ws.Range("V7").Value = [double]"0.03483694568883143"

# row 8
$ws.Range("B8").Value = [double]"0.0005142141005934914"
$ws.Range("C8").Value = [double]"0.1468234750174995"
$ws.Range("D8").Value = [double]"0.004967089149840875"
$ws.Range("E8").Value = [double]"0.006680156260488962"
$ws.Range("F8").Value = [double]"0.4259403133055403"
$ws.Range("G8").Value = [double]"0.08101698208793137"
$ws.Range("H8").Value = [double]"4.241841841375875E-05"
$ws.Range("I8").Value = [double]"6.119401243570504E-06"
$ws.Range("J8").Value = [double]"0.006936706887952397"
$ws.Range("K8").Value = [double]"0.1224556766298276"
$ws.Range("L8").Value = [double]"0.007252575716014791"
$ws.Range("M8").Value = [double]"0.01530564068213586"
$ws.Range("N8").Value = [double]"0.002288302281484026"
$ws.Range("O8").Value = [double]"3.214693312517258E-08"
$ws.Range("P8").Value = [double]"0.0002243250443006807"
$ws.Range("Q8").Value = [double]"0.0004438420140741797"
$ws.Range("R8").Value = [double]"0.0004785348084436246"
$ws.Range("S8").Value = [double]"0.001755203861496419"
$ws.Range("T8").Value = [double]"7.551334526633466E-10"
$ws.Range("U8").Value = [double]"3.325322783403728E-07"
$ws.Range("V8").Value = [double]"0.03133095944163506"

# row 9
$ws.Range("B9").Value = [double]"6.700317159180644E-06"
$ws.Range("C9").Value = [double]"0.002030465225620723"
$ws.Range("D9").Value = [double]"6.742453047929233E-05"
$ws.Range("E9").Value = [double]"8.951435836395369E-05"
$ws.Range("F9").Value = [double]"0.006561235673855828"
$ws.Range("G9").Value = [double]"0.0009081999569219748"
$ws.Range("H9").Value = [double]"7.40766342571984E-07"
$ws.Range("I9").Value = [double]"2.78930871693052E-07"
$ws.Range("J9").Value = [double]"8.986617930536877E-05"
$ws.Range("K9").Value = [double]"0.001906933140010266"
$ws.Range("L9").Value = [double]"0.000154917789506855"
$ws.Range("M9").Value = [double]"0.0007143758247708101"
$ws.Range("N9").Value = [double]"0.0006225583469750111"
$ws.Range("O9").Value = [double]"8.261620183242931E-10"
$ws.Range("P9").Value = [double]"3.365580160921686E-06"
$ws.Range("Q9").Value = [double]"5.30122366023758E-06"
$ws.Range("R9").Value = [double]"5.624321345617954E-06"
$ws.Range("S9").Value = [double]"3.093410026360511E-05"
$ws.Range("T9").Value = [double]"1.47954789967821E-11"
$ws.Range("U9").Value = [double]"4.807741455253776E-09"
$ws.Range("V9").Value = [double]"0.0004624424323723597"

# row 10
$ws.Range("B10").Value = [double]"9.422588095789063E-05"
$ws.Range("C10").Value = [double]"0.03243189536900044"
$ws.Range("D10").Value = [double]"0.00110678600668379"
$ws.Range("E10").Value = [double]"0.001427522988424535"
$ws.Range("F10").Value = [double]"0.1025072573169607"
$ws.Range("G10").Value = [double]"0.01432989096197988"
$ws.Range("H10").Value = [double]"9.109993019665222E-06"
$ws.Range("I10").Value = [double]"1.073992531681E-06"
$ws.Range("J10").Value = [double]"0.001538819993760139"
$ws.Range("K10").Value = [double]"0.02373957277752879"
$ws.Range("L10").Value = [double]"0.00205120860092341"
$ws.Range("M10").Value = [double]"0.0007037848907318007"
$ws.Range("N10").Value = [double]"0.000554922257555286"
$ws.Range("O10").Value = [double]"8.269592878411739E-09"
$ws.Range("P10").Value = [double]"4.537142469902325E-05"
$ws.Range("Q10").Value = [double]"7.16094306449583E-05"
$ws.Range("R10").Value = [double]"7.6701961563405E-05"
$ws.Range("S10").Value = [double]"0.0003679586712385384"
$ws.Range("T10").Value = [double]"1.35140213923397E-10"
$ws.Range("U10").Value = [double]"7.003126439395945E-08"
$ws.Range("V10").Value = [double]"0.005007019334926147"

# row 11
$ws.Range("B11").Value = [double]"-0.00204419637904557"
$ws.Range("C11").Value = [double]"-2.481352495021464"
$ws.Range("D11").Value = [double]"-0.04186986718212457"
$ws.Range("E11").Value = [double]"-0.05713111537631673"
$ws.Range("F11").Value = [double]"-2.287035033970398"
$ws.Range("G11").Value = [double]"-0.2480134799510296"
$ws.Range("H11").Value = [double]"-0.0006770049331384043"
$ws.Range("I11").Value = [double]"-5.377657245532569E-05"
$ws.Range("J11").Value = [double]"-0.2721480275830908"
$ws.Range("K11").Value = [double]"-0.9818855493608475"
$ws.Range("L11").Value = [double]"-0.03679645533388078"
$ws.Range("M11").Value = [double]"-0.009861756353170511"
$ws.Range("N11").Value = [double]"-0.03895400170292455"
$ws.Range("O11").Value = [double]"-6.581680745704826E-07"
$ws.Range("P11").Value = [double]"-0.0009343549484718144"
$ws.Range("Q11").Value = [double]"-0.001837736491969356"
$ws.Range("R11").Value = [double]"-0.001970886450070403"
$ws.Range("S11").Value = [double]"-0.008419273731913845"
$ws.Range("T11").Value = [double]"-8.238229938947979E-09"
$ws.Range("U11").Value = [double]"-4.038304877756715E-06"
$ws.Range("V11").Value = [double]"-0.1021460400779775"

# row 12
$ws.Range("B12").Value = [double]"-0.002044196379189382"
$ws.Range("C12").Value = [double]"-2.481352494593619"
$ws.Range("D12").Value = [double]"-0.04186986702865887"
$ws.Range("E12").Value = [double]"-0.05713111517860117"
$ws.Range("F12").Value = [double]"-2.287035021317406"
$ws.Range("G12").Value = [double]"-0.2480134797703242"
$ws.Range("H12").Value = [double]"-0.0006770049328409821"
$ws.Range("I12").Value = [double]"-5.377657244069622E-05"
$ws.Range("J12").Value = [double]"-0.2721480274703113"
$ws.Range("K12").Value = [double]"-0.9818855471327521"
$ws.Range("L12").Value = [double]"-0.03679645534047439"
$ws.Range("M12").Value = [double]"-0.009861756715659078"
$ws.Range("N12").Value = [double]"-0.03895400169068426"
$ws.Range("O12").Value = [double]"-6.581680747703665E-07"
$ws.Range("P12").Value = [double]"-0.0009343549490316631"
$ws.Range("Q12").Value = [double]"-0.001837736497052592"
$ws.Range("R12").Value = [double]"-0.00197088645468444"
$ws.Range("S12").Value = [double]"-0.008419273728277915"
$ws.Range("T12").Value = [double]"-8.238229941122425E-09"
$ws.Range("U12").Value = [double]"-4.038304876826917E-06"
$ws.Range("V12").Value = [double]"-0.1021460400133375"
